$d = $word.ActiveDocument

# 1) Update the release/revision metadata lines at the top of the document.
$d.Content.Find.Execute("Release name: 2.0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Release name: 1.0", 2) | Out-Null

$d.Content.Find.Execute("Release date: February 21, 2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Release date: March 11, 2018", 2) | Out-Null

$d.Content.Find.Execute("Revision number: 2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Revision number: 3", 2) | Out-Null

$d.Content.Find.Execute("Revision date: February 7, 2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Revision date: March 11, 2018", 2) | Out-Null

# 2) Remove the "row#(player) constraint" user story row from the first table
#    (it currently sits between the "savable configurations" row and the
#    "appealing application" row).
$t1 = $d.Tables.Item(1)
for ($i = 1; $i -le $t1.Rows.Count; $i++) {
    $row = $t1.Rows.Item($i)
    $text = $row.Cells.Item(1).Range.Text
    if ($text -like "*row#(player) is selected so I can consolidate my risks*") {
        $row.Delete()
        break
    }
}

# 3) Re-add that same user story, in simplified form, as a brand new last row
#    of the second (difficulty) table.
$t2 = $d.Tables.Item(2)
$newRow = $t2.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "As a DFS player I want a ‘constraint’ for the number of times a certain row#(player) is selected so I can consolidate my risks"
$newRow.Cells.Item(2).Range.Text = "5"
